# Insert a new data row at row 492 (pushes the existing rows 492..597 down
# to 493..598, growing the sheet's used range from A1:R597 to A1:R598) and
# populate it with the new weekly price-report record for Acelga.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(492).Insert()

$ws.Range("A492").Value = 10
$ws.Range("B492").Value = "Vega Modelo de Temuco"
$ws.Range("C492").Value = "La Araucanía"
$ws.Range("D492").Value = 45244
$ws.Range("E492").Value = 9
$ws.Range("F492").Value = 100112009
$ws.Range("G492").Value = "Acelga"
$ws.Range("H492").Value = "Sin especificar"
$ws.Range("I492").Value = "Primera"
$ws.Range("J492").Value = 65
$ws.Range("K492").Value = 10000
$ws.Range("L492").Value = 10000
$ws.Range("M492").Value = 10000
$ws.Range("N492").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O492").Value = "Provincia de Cautín"
$ws.Range("P492").Value = 833
$ws.Range("Q492").Value = 12
$ws.Range("R492").Value = "Hortaliza"
